$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - refreshed crypto symbol data.
# Values are stored as text in the sheet, so we use a leading apostrophe
# to force Excel to keep them as text instead of auto-converting to numbers.
$ws.Range("D2").Value = "'271.73"
$ws.Range("D3").Value = "'23.09"
$ws.Range("D4").Value = "'6.358"
$ws.Range("D5").Value = "'0.06322"
$ws.Range("D6").Value = "'3.661"
$ws.Range("D7").Value = "'6.771"
$ws.Range("D8").Value = "'1.394"
$ws.Range("D9").Value = "'0.8371"
$ws.Range("D10").Value = "'0.1626"
$ws.Range("D11").Value = "'0.08362"
$ws.Range("D12").Value = "'0.03432"
$ws.Range("D13").Value = "'0.03154"
$ws.Range("D14").Value = "'0.09318"
$ws.Range("D15").Value = "'3.915"
$ws.Range("D16").Value = "'0.001709"
$ws.Range("D17").Value = "'0.04858"
$ws.Range("D18").Value = "'0.006270"
$ws.Range("D19").Value = "'0.005499"
$ws.Range("D20").Value = "'0.001087"
$ws.Range("D21").Value = "'0.0001496"
$ws.Range("D22").Value = "'3.739"
$ws.Range("D23").Value = "'2.343"
$ws.Range("D24").Value = "'0.01386"
$ws.Range("D25").Value = "'0.3382"
$ws.Range("D27").Value = "'0.0002675"
$ws.Range("D40").Value = "'0.04690"
$ws.Range("D41").Value = "'0.006892"
$ws.Range("D42").Value = "'0.1179"
$ws.Range("D43").Value = "'0.003449"
$ws.Range("D44").Value = "'0.01255"
$ws.Range("D45").Value = "'0.00006256"
$ws.Range("D46").Value = "'0.00000000748"
$ws.Range("D47").Value = "'0.6980"
$ws.Range("D48").Value = "'0.1219"
$ws.Range("D49").Value = "'0.00002094"
$ws.Range("D50").Value = "'0.01236"

Write-Host "Updated price column for 36 rows"
